$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.594.71'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.39%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.219.44'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.86%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '240.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.81%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.13%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '74.83'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.07%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('E9').Value = '  +1.43%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.20'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.22%  '

$ws.Range('E11').Value = '  -2.18%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.99'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.66%  '

$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('E14').Value = '  -2.07%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.551.77'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.88%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.65'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.33%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.214.28'
$ws.Range('D17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.800'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.12%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '42.501.04'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.47%  '

$ws.Range('E20').Value = '  +0.53%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '70.83'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.09%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.91'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.55%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.06'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -8.76%  '

$ws.Range('E24').Value = '  -0.16%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.15'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +5.61%  '

$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.90'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.68%  '

$ws.Range('E28').Value = '  -6.84%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.22'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.21%  '

$ws.Range('E30').Value = '  -0.65%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.51'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.23%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '34.65'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +12.73%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.24'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.47%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0790'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.82%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.32'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.29%  '

$ws.Range('E36').Value = '  -1.92%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.107'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.16%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.37'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.30%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0320'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.44%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '12.60'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.29%  '

$ws.Range('E41').Value = '  +0.49%  '

$ws.Range('E42').Value = '  -2.88%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '60.62'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.64%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.197'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.23%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.54'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.60%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0985'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.56%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '99.40'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.34%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.10'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.24%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.31'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.58%  '

$ws.Range('E50').Value = '  -2.86%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.428'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +16.20%  '
